$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell = 'D2'; Value = '44.208.50'},
    @{Cell = 'E2'; Value = '  +0.20%  '},
    @{Cell = 'D3'; Value = '2.242.14'},
    @{Cell = 'E3'; Value = '  -0.08%  '},
    @{Cell = 'E4'; Value = '  +0.21%  '},
    @{Cell = 'D5'; Value = '306.89'},
    @{Cell = 'E5'; Value = '  -2.70%  '},
    @{Cell = 'D6'; Value = '94.62'},
    @{Cell = 'E6'; Value = '  -4.72%  '},
    @{Cell = 'E7'; Value = '  -0.51%  '},
    @{Cell = 'D8'; Value = '1.00'},
    @{Cell = 'E8'; Value = '  +0.28%  '},
    @{Cell = 'E9'; Value = '  -1.61%  '},
    @{Cell = 'D10'; Value = '34.70'},
    @{Cell = 'E10'; Value = '  -4.04%  '},
    @{Cell = 'D11'; Value = '0.0811'},
    @{Cell = 'E11'; Value = '  -1.54%  '},
    @{Cell = 'E12'; Value = '  -2.53%  '},
    @{Cell = 'D14'; Value = '2.584.05'},
    @{Cell = 'E14'; Value = '  +0.07%  '},
    @{Cell = 'D15'; Value = '2.248.57'},
    @{Cell = 'E15'; Value = '  +0.57%  '},
    @{Cell = 'E16'; Value = '  -1.59%  '},
    @{Cell = 'E17'; Value = '  -3.23%  '},
    @{Cell = 'D18'; Value = '43.977.29'},
    @{Cell = 'E18'; Value = '  +0.19%  '},
    @{Cell = 'D19'; Value = '0.0₃0963'},
    @{Cell = 'E19'; Value = '  -1.24%  '},
    @{Cell = 'D20'; Value = '6.38'},
    @{Cell = 'E20'; Value = '  +0.93%  '},
    @{Cell = 'D21'; Value = '12.11'},
    @{Cell = 'E21'; Value = '  -8.12%  '},
    @{Cell = 'D22'; Value = '65.56'},
    @{Cell = 'E22'; Value = '  -0.62%  '},
    @{Cell = 'D23'; Value = '237.71'},
    @{Cell = 'E23'; Value = '  -0.22%  '},
    @{Cell = 'E24'; Value = '  -0.81%  '},
    @{Cell = 'E25'; Value = '  -1.03%  '},
    @{Cell = 'E26'; Value = '  -0.09%  '},
    @{Cell = 'B27'; Value = 'InjectiveProtocol'},
    @{Cell = 'C27'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'},
    @{Cell = 'D27'; Value = '38.30'},
    @{Cell = 'E27'; Value = '  +5.25%  '},
    @{Cell = 'B28'; Value = 'Cosmos'},
    @{Cell = 'C28'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'},
    @{Cell = 'D28'; Value = '9.92'},
    @{Cell = 'E28'; Value = '  -2.50%  '},
    @{Cell = 'D29'; Value = '2.21'},
    @{Cell = 'E29'; Value = '  +3.65%  '},
    @{Cell = 'D30'; Value = '20.02'},
    @{Cell = 'E30'; Value = '  -0.23%  '},
    @{Cell = 'D31'; Value = '5.84'},
    @{Cell = 'E31'; Value = '  -2.50%  '},
    @{Cell = 'D32'; Value = '153.16'},
    @{Cell = 'E32'; Value = '  -0.97%  '},
    @{Cell = 'E33'; Value = '  -4.65%  '},
    @{Cell = 'E34'; Value = '  -1.55%  '},
    @{Cell = 'D35'; Value = '3.15'},
    @{Cell = 'E35'; Value = '  -5.67%  '},
    @{Cell = 'E37'; Value = '  -0.95%  '},
    @{Cell = 'D38'; Value = '1.77'},
    @{Cell = 'E38'; Value = '  -7.48%  '},
    @{Cell = 'D39'; Value = '3.53'},
    @{Cell = 'E39'; Value = '  +0.18%  '},
    @{Cell = 'E40'; Value = '  -4.93%  '},
    @{Cell = 'E41'; Value = '  -9.17%  '},
    @{Cell = 'E42'; Value = '  -2.93%  '},
    @{Cell = 'E43'; Value = '  +0.22%  '},
    @{Cell = 'D44'; Value = '1.744.14'},
    @{Cell = 'E44'; Value = '  +2.44%  '},
    @{Cell = 'D45'; Value = '82.48'},
    @{Cell = 'E45'; Value = '  +0.25%  '},
    @{Cell = 'E46'; Value = '  -2.40%  '},
    @{Cell = 'B47'; Value = 'EnergySwap'},
    @{Cell = 'C47'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'},
    @{Cell = 'D47'; Value = '14.97'},
    @{Cell = 'E47'; Value = '  +8.52%  '},
    @{Cell = 'B48'; Value = 'Aave'},
    @{Cell = 'C48'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'},
    @{Cell = 'D48'; Value = '99.79'},
    @{Cell = 'E48'; Value = '  -1.88%  '},
    @{Cell = 'B49'; Value = 'THORChain'},
    @{Cell = 'C49'; Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'},
    @{Cell = 'D49'; Value = '4.93'},
    @{Cell = 'E49'; Value = '  -4.82%  '},
    @{Cell = 'B50'; Value = 'FraxShare'},
    @{Cell = 'C50'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'},
    @{Cell = 'D50'; Value = '8.06'},
    @{Cell = 'E50'; Value = '  -1.00%  '},
    @{Cell = 'B51'; Value = 'Stacks'},
    @{Cell = 'C51'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'},
    @{Cell = 'D51'; Value = '1.57'},
    @{Cell = 'E51'; Value = '  -3.89%  '},
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
